$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.132.97"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "3.382.39"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.95"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "658.26"
$ws.Range("E6").Value = "  -1.87%  "
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.423"
$ws.Range("E8").Value = "  -3.68%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  -2.49%  "
$ws.Range("D11").Value = "3.379.82"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("E12").Value = "  -3.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.92"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "97.765.03"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.09"
$ws.Range("E15").Value = "  -5.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000256"
$ws.Range("E16").Value = "  -4.48%  "
$ws.Range("D17").Value = "4.017.04"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.16"
$ws.Range("E18").Value = "  +2.62%  "
$ws.Range("D19").Value = "3.379.51"
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.97"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.519"
$ws.Range("E21").Value = "  -8.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.07"
$ws.Range("E22").Value = "  -1.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "509.45"
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("E24").Value = "  -1.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000201"
$ws.Range("E25").Value = "  -3.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.87"
$ws.Range("E26").Value = "  +2.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.50"
$ws.Range("E27").Value = "  -5.01%  "
$ws.Range("E28").Value = "  -5.07%  "
$ws.Range("D29").Value = "3.567.22"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.62"
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("E31").Value = "  -4.76%  "
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.188"
$ws.Range("E33").Value = "  -5.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.60"
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "28.72"
$ws.Range("E37").Value = "  -4.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.89"
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("E39").Value = "  -5.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "526.75"
$ws.Range("E40").Value = "  -3.06%  "
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.844"
$ws.Range("E44").Value = "  -4.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.72"
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0424"
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.69"
$ws.Range("E47").Value = "  -3.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.59"
$ws.Range("E48").Value = "  -4.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.23"
$ws.Range("E49").Value = "  +4.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "55.41"
$ws.Range("E50").Value = "  +2.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.53"
$ws.Range("E51").Value = "  -5.49%  "
